$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- text blocks (single-quoted here-strings: literal, no $var expansion) ---
$b119 = @'
Collection: Aggregate Operations 
'@
$b120 = @'
Lambda Expression
'@
$c120 = @'
# What is it 
Lambda expression looks like a method declaration, you can consider lambda expressions as anonymous methods- methods without a name.
# Where can use
## To replace any anonymous inner class that has only **single method interface**. In below example, tester1 and tester2 are equivalent.
```
interface Foo{
    boolean check(T t);
}
public void main(){
    tester1.setFoo(new Foo{
        @Override
        public boolean check(T t){
            return t.foobar();
        }
    });
    tester2.setFoo( t -> t.foobar() );
}
```
## To pass into **Aggregate Operations** that accept Lambda Expressions as parameters
```
roster
    .stream()
    .filter(p -> p.getGender() == Person.Sex.MALE && p.getAge() >= 18 && p.getAge() <= 25)
    .map(p -> p.getEmailAddress())
    .forEach(email -> System.out.println(email));
```
# Syntax
For below code, the statement after -> will treat as a return statement automatically:
` p -> p.getAge() `
For below code, return is implicitly decare:
```
p -> { 
    p.refresh();
    p.increment();
    return p.getAge();
}
```
'@
$b121 = @'
Collection: Impls
'@
$c121 = @'
# Implementations of Collection Interface
* Impl only affect performance, not feature
* Vector and Hashtable are legacy collections, they are synchonized. But lastest alternatives ArrayList and HashMap are not, so to reduce sync overhead and unnecessary feature. " In general, it is good API design practice not to make users pay for a feature they don't use."
* If thread-safe collections needed:
    * Collection impl from java.util.concurrent (offers much higher concurrency than sync wrapping collection)
    * Synchronization wrapper to wrap your collections
# Commonly use general purpose collection implementations
* For the Set interface, HashSet is the most commonly used implementation.
* For the List interface, ArrayList is the most commonly used implementation.
* For the Map interface, HashMap is the most commonly used implementation.
* For the Queue interface, LinkedList is the most commonly used implementation.
* For the Deque interface, ArrayDeque is the most commonly used implementation.
'@
$b122 = @'
Collection: Algorithms
'@
$c122 = @'
The collection interface does not carry operation itself (sort, search, shuffle and etc). These operation and algorithms are provided in Collections:
Collections.min(List, Comparator)
Collections.max(List, Comparator)
Collections.sort(List, Comparator)
Collections.binarySearch(List, T)
'@
$b123 = @'
Collection: Custom Collection Design
'@
$c123 = @'
AbstractList is provided to resue for new Impl of Collections. As long as user impl the constructor, the get, the set and the size methods, the other bulk operations are provided in AbstractList already.
For other interface, the same principle applies.
'@
$b124 = @'
Collection: Inter-operability
'@
$c124 = @'
# Compatibility 
* Upward Compatibility: Get array from oldMethod() and pass to newMethod with Arrays.asList() that construct an List based on the array
* Backward Compatibility: Get collection from newMethod() and "down-grade" it to array with List.toArray() method and pass the arry to oldMethod()
* I learnt that 兼容性在不同語境下有不同含義. 在JVM層面下, 那就是JRE7可否運行JDK6的程序這樣的問題. 而在代碼層面, 則是接口過渡問題. 如上所述, 只要能夠將舊接口的訊息不失真地轉移到新接口, 則稱其Upward Compability; 而新接口的東西, 反之能轉換回舊接口可識別的訊息, 則稱Backward Compatibility
'@

# Row 119 ("RTFM" / "Aggregate Operations " / <intro text>) is the template:
# clone its formatting (style, row height) into 5 fresh rows below it by
# copy+insert, once per row (the clipboard is consumed by each Insert, so we
# re-copy before every iteration). This keeps the new rows' cell style in
# sync with every other entry in the RTFM table instead of picking up a
# blank default style.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(119).Copy()
    $ws.Rows.Item(120).Insert(-4121, -4163)
}
# The insert pushed the old (mostly-empty) row 120 - which only ever held
# "A120" = RTFM with no B/C - down to row 125; drop that leftover duplicate.
$ws.Rows.Item(125).Delete() | Out-Null

# Fill row 120 first (new "Lambda Expression" entry) before renaming row 119's
# title below, so new shared-string entries land in the same order as the
# target edit.
$ws.Range("B120").Value = $b120
$ws.Range("C120").Value = $c120

# Row 119: retitle "Aggregate Operations " -> "Collection: Aggregate Operations ".
# Its content cell (C119) is left untouched - it was already correct from the
# row copy above.
$ws.Range("B119").Value = $b119

# Rows 121-124: new Collection entries (A121:A124 already hold "RTFM" from the
# row copy above).
$ws.Range("B121").Value = $b121
$ws.Range("C121").Value = $c121

$ws.Range("B122").Value = $b122
$ws.Range("C122").Value = $c122

$ws.Range("B123").Value = $b123
$ws.Range("C123").Value = $c123

$ws.Range("B124").Value = $b124
$ws.Range("C124").Value = $c124

# Pin every touched row back to the table's standard height (large pasted text
# would otherwise trigger autofit to a taller row).
for ($r = 119; $r -le 124; $r++) {
    $ws.Rows.Item($r).RowHeight = 32.25
}

# Leave the same cell selected as in the target workbook.
$ws.Range("C124").Select() | Out-Null
